$wb = $excel.ActiveWorkbook

# Rename "Hoja1" -> "cables"
$wsCables = $wb.Worksheets.Item("Hoja1")
$wsCables.Name = "cables"

# Add new "fusibles" sheet right after "cables", before "Hoja2"
$wsHoja2 = $wb.Worksheets.Item("Hoja2")
$wsFusibles = $wb.Worksheets.Add($null, $wsCables)
$wsFusibles.Name = "fusibles"

# Populate the "fusibles" sheet
$wsFusibles.Range("A1").Value = "Fusibles (A)"
$wsFusibles.Range("A2").Value = 50
$wsFusibles.Range("A3").Value = 80
$wsFusibles.Range("A4").Value = 100
$wsFusibles.Range("A5").Value = 125
$wsFusibles.Range("A6").Value = 200
$wsFusibles.Range("A7").Value = 250
$wsFusibles.Range("A8").Value = 300

# Update selection on Hoja2 sheet, then restore "cables" as the active tab
$wsHoja2 = $wb.Worksheets.Item("Hoja2")
$wsHoja2.Range("B2").Select()
$wsCables.Activate()

$wb.Save()
